$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.297.48"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.81%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.704.56"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "223.76"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5303"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  -0.13%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2654"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06569"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.68%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.73"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.47%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07630"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.513"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.05%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.940.83"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.705.03"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.5765"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0₅8125"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "67.51"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "27.308.82"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.84%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "215.04"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("E20").Value = "  -0.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.611"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.08%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.36"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.951"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  -2.95%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.705"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.1201"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.205"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("E29").Value = "  -4.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.05367"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.86%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.287"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.42%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.463"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.49%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.399"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.10%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.636"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.870"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.76%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.414"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9440"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.26%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.5785"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.77%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01627"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.60%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.757"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.037.86"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.11%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8391"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.41%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "100.90"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.57%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.848.69"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  +1.42%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "57.64"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  +1.81%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.070"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.09%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05229"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
